$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()
$ws.Range("T13").Value = 2
$ws.Range("T15").Select()
